$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.811.30"
$ws.Range("E2").Value = "  +4.94%  "
$ws.Range("D3").Value = "2.658.48"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'568.90"
$ws.Range("E5").Value = "  +7.00%  "
$ws.Range("D6").Value = "'146.32"
$ws.Range("E6").Value = "  +3.90%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.609"
$ws.Range("E8").Value = "  +4.60%  "
$ws.Range("D9").Value = "2.658.80"
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +6.24%  "
$ws.Range("E12").Value = "  +7.29%  "
$ws.Range("D13").Value = "'0.344"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").Value = "3.104.91"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "60.722.87"
$ws.Range("E15").Value = "  +4.90%  "
$ws.Range("D16").Value = "'22.10"
$ws.Range("E16").Value = "  +7.53%  "
$ws.Range("D17").Value = "'0.0000139"
$ws.Range("E17").Value = "  +6.37%  "
$ws.Range("D18").Value = "2.657.34"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").Value = "'4.55"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("D20").Value = "'343.14"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").Value = "'10.46"
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("E22").Value = "  +3.94%  "
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'66.31"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "'0.439"
$ws.Range("E25").Value = "  +5.46%  "
$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").Value = "  +6.22%  "
$ws.Range("D29").Value = "0.0₃0817"
$ws.Range("E29").Value = "  +13.57%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +5.99%  "
$ws.Range("E32").Value = "  +6.63%  "
$ws.Range("D33").Value = "'160.09"
$ws.Range("E33").Value = "  +3.77%  "
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("E35").Value = "  +6.93%  "
$ws.Range("D36").Value = "'0.901"
$ws.Range("E36").Value = "  +9.66%  "
$ws.Range("E37").Value = "  +6.66%  "
$ws.Range("E38").Value = "  +10.53%  "
$ws.Range("E39").Value = "  +8.56%  "
$ws.Range("D40").Value = "'37.46"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "'300.39"
$ws.Range("E41").Value = "  +8.51%  "
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'0.0986"
$ws.Range("E44").Value = "  +4.97%  "
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("D46").Value = "'0.0546"
$ws.Range("E46").Value = "  +4.15%  "
$ws.Range("D47").Value = "'128.41"
$ws.Range("E47").Value = "  +15.46%  "
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").Value = "'0.0236"
$ws.Range("E50").Value = "  +5.49%  "
$ws.Range("D51").Value = "'4.66"
$ws.Range("E51").Value = "  +6.60%  "
